$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correção dos rótulos da linha de cabeçalho 2: as células com textos
# genéricos "unnamed: 1_level_1" e "unnamed: 5_level_1" (gerados pelo
# pandas) são corrigidas para "total", conforme os demais blocos da tabela.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
